# Weekly price-list refresh: a new "Arveja Verde" record for
# Mercado Mayorista Lo Valledor de Santiago is inserted at row 206,
# pushing every subsequent record down by one row (old row 292 becomes
# the new last row, 293).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record by inserting a blank row at 206;
# everything that used to live at row 206 (and below) shifts to row 207+.
$ws.Rows.Item(206).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(206, 1).Value  = 6
$ws.Cells.Item(206, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(206, 3).Value  = "Metropolitana"
$ws.Cells.Item(206, 4).Value  = 45009
$ws.Cells.Item(206, 5).Value  = 13
$ws.Cells.Item(206, 6).Value  = 100112022
$ws.Cells.Item(206, 7).Value  = "Arveja Verde"
$ws.Cells.Item(206, 8).Value  = "Perfection"
$ws.Cells.Item(206, 9).Value  = "Primera"
$ws.Cells.Item(206, 10).Value = 800
$ws.Cells.Item(206, 11).Value = 25000
$ws.Cells.Item(206, 12).Value = 27000
$ws.Cells.Item(206, 13).Value = 26125
$ws.Cells.Item(206, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(206, 15).Value = "Región Metropolitana"
$ws.Cells.Item(206, 16).Value = 1045
$ws.Cells.Item(206, 17).Value = 25
$ws.Cells.Item(206, 18).Value = "Hortaliza"
